$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 gains two new notes in D3 and F3, written in the same
# "Times New Roman, 12pt" body font already used elsewhere on the sheet
# (e.g. C2:F2), but without the wrap/top alignment those header cells use.
# Borrow the font via a format-only paste from C2 (so the existing font
# record is reused instead of a new one being minted), then turn wrapping
# back off before filling in the actual text.
$ws.Range("C2").Copy() | Out-Null
$ws.Range("D3").PasteSpecial(-4122) | Out-Null
$ws.Range("D3").WrapText = $false
$ws.Range("D3").Value = "Ogarnia muzke "

$ws.Range("C2").Copy() | Out-Null
$ws.Range("F3").PasteSpecial(-4122) | Out-Null
$ws.Range("F3").WrapText = $false
$ws.Range("F3").Value = "Niedziela dzień wolny 😴"

# Clear the marching-ants marquee left behind by the copy operations.
$excel.CutCopyMode = $false

# The active selection moved from N10 to F4.
$ws.Range("F4").Select()
